$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 (I0) and J1 (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style of the existing header row (bold/border/centered), copied from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for I2:J51 (I-value, J-value) per row
$data = @(
@(6,6),@(4,5),@(10,10),@(7,8),@(6,7),@(6,6),@(6,6),@(9,9),@(9,10),@(7,7),
@(7,7),@(7,7),@(8,8),@(6,7),@(8,8),@(6,7),@(8,8),@(7,7),@(8,8),@(10,11),
@(4,5),@(9,9),@(6,7),@(7,7),@(7,7),@(6,6),@(9,9),@(6,6),@(7,8),@(5,6),
@(9,9),@(7,7),@(8,8),@(8,8),@(9,9),@(7,7),@(5,6),@(7,8),@(6,6),@(5,5),
@(7,8),@(5,6),@(7,8),@(7,7),@(5,5),@(6,6),@(6,6),@(7,7),@(9,9),@(7,7)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $data[$idx][0]
    $ws.Cells.Item($r, 10).Value = $data[$idx][1]
}
